$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-15 down to 4-16
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 with the new "Agrícola del Norte S.A. de Arica" record
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44764
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 100112017
$ws.Range("G3").Value = "Ramas de apio"
$ws.Range("H3").Value = "Americana (o)"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 7500
$ws.Range("N3").Value = "`$/atado 7 kilos"
$ws.Range("O3").Value = "Región de Arica y Parinacota"
$ws.Range("P3").Value = 7500
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
